# Applies the commit "custom exception, decode password":
#  - RUNMANAGER: toggle the "execute" flag for loginLogoutTest (C2) and newTest (C4)
#  - DATA: change one browser value, toggle a couple of "execute" flags, and
#          replace the plaintext password "admin123" with its base64-encoded
#          form "YWRtaW4xMjM=" wherever it appears in the password column
#  - both sheets: move the active selection

$wb = $excel.ActiveWorkbook

$wsRunManager = $wb.Worksheets.Item("RUNMANAGER")
$wsData = $wb.Worksheets.Item("DATA")

# --- RUNMANAGER sheet -------------------------------------------------
$wsRunManager.Range("C2").Value = "yes"
$wsRunManager.Range("C4").Value = "no"

$wsRunManager.Range("C3").Select()

# --- DATA sheet ---------------------------------------------------------
$wsData.Range("C2").Value = "chrome"

$wsData.Range("E2").Value = "YWRtaW4xMjM="
$wsData.Range("E3").Value = "YWRtaW4xMjM="
$wsData.Range("E4").Value = "YWRtaW4xMjM="
$wsData.Range("E5").Value = "YWRtaW4xMjM="
$wsData.Range("E6").Value = "YWRtaW4xMjM="

$wsData.Range("B4").Value = "no"
$wsData.Range("B7").Value = "no"
$wsData.Range("B8").Value = "no"

$wsData.Range("C2").Select()
